$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = -0.6762232433871019
$ws.Range("J5").Value = 0.451917121526646
$ws.Range("K5").Value = 0.2157933155075947
$ws.Range("L5").Value = 2.611954445245234
